# Update project plan sheet: fill in completion percentages for row 3-7,
# add a remark for row 7, update the summary text for the first week block,
# and move the active selection to G9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Completion percentage column (C) for the first block of tasks (rows 3-6 = 100%, row 7 = 60%)
$percentRange = $ws.Range("C3:C7")
$percentRange.NumberFormat = "0%"

$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 0.6

# Update the weekly summary text (merged cell A8:D9)
$ws.Range("A8").Value = "总结：已进行项目需求分析讨论，并通过分析确定每个功能模块的功能"

# Remark for the partially completed task
$ws.Range("D7").Value = "暂时未得知如何设计用例图"

# Move the selection/active cell to G9, matching the saved view state
$ws.Range("G9").Select() | Out-Null
